$d = $word.ActiveDocument

# 1. Title font size: 40 half-points (20pt) -> 32 half-points (16pt)
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.MoveEnd(1, -1) | Out-Null
$titleRange.Font.Size = 16

# 2. Replace the five empty paragraphs with a paragraph containing a long
#    underscore divider line. Empty paragraphs are located (by index, from
#    the original document) right after: the contact-info paragraph, the
#    education paragraph, the Coresoft experience paragraph, the perfil
#    profesional paragraph, and the Clothing Ecommerce project paragraph.
$divider = "________________________________________________________________________________"

$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "`r") {
        $p.Range.Text = $divider
    }
}

# 3. Update the Coresoft Front-End experience description text.
$old1 = "Formé parte del equipo de desarrollo Front-End, contribuyendo a la creación de interfaces web y móviles modernas, optimizadas para rendimiento y usabilidad, asegurando buenas prácticas de diseño y desarrollo."
$new1 = "Formé parte del equipo de desarrollo Front-End, contribuyendo a la creación de interfaces web y móviles optimizadas, garantizando un rendimiento eficiente y una experiencia de usuario fluida."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# 4. Update the "Normal" style's default font from Calibri to Georgia.
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Name = "Georgia"
